# Natmi following Dr Hou advice
# Re-run of the NATMI LR-pair export for Tnc -> Egfr now includes a 4th
# sending/target cluster ("M2", in addition to ECs, FAPs and sCs) and
# refreshed statistics for every sending x target cluster combination.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one data row: A..T in column order
# (Sending cluster, Ligand symbol, Receptor symbol, Target cluster, then
# the 16 numeric NATMI statistic columns).
$data = @(
    @("ECs","Tnc","Egfr","ECs",3,1,8.948174,26.844522,0.0695931738232498,0.0695931738232498,3,1,1.854571666666667,5.563715,0.01651371646154392,0.01651371646154392,16.59502996880333,149.35526971923,0.001149241940176087,0.001149241940176088),
    @("ECs","Tnc","Egfr","FAPs",3,1,8.948174,26.844522,0.0695931738232498,0.0695931738232498,3,1,92.91372433333333,278.741173,0.8273343794712995,0.8273343794712996,831.4081723227007,7482.673550904306,0.05757682528049665,0.05757682528049666),
    @("ECs","Tnc","Egfr","M2",3,1,8.948174,26.844522,0.0695931738232498,0.0695931738232498,3,1,0.162136,0.4864080000000001,0.001443712303133186,0.001443712303133187,1.450821139664,13.057390256976,0.0001004725212627121,0.0001004725212627121),
    @("ECs","Tnc","Egfr","sCs",3,1,8.948174,26.844522,0.0695931738232498,0.0695931738232498,3,1,17.37449166666667,52.123475,0.1547081917640233,0.1547081917640233,155.4699745948833,1399.22977135395,0.01076663408131434,0.01076663408131434),
    @("FAPs","Tnc","Egfr","ECs",3,1,105.8801323333333,317.640397,0.8234679448457706,0.8234679448457706,3,1,1.854571666666667,5.563715,0.01651371646154392,0.01651371646154392,196.3622934883172,1767.260641394855,0.01359851615635334,0.01359851615635334),
    @("FAPs","Tnc","Egfr","FAPs",3,1,105.8801323333333,317.640397,0.8234679448457706,0.8234679448457706,3,1,92.91372433333333,278.741173,0.8273343794712995,0.8273343794712996,9837.717427996187,88539.45685196569,0.6812833411634819,0.681283341163482),
    @("FAPs","Tnc","Egfr","M2",3,1,105.8801323333333,317.640397,0.8234679448457706,0.8234679448457706,3,1,0.162136,0.4864080000000001,0.001443712303133186,0.001443712303133187,17.16698113599734,154.502830223976,0.001188850803209639,0.001188850803209639),
    @("FAPs","Tnc","Egfr","sCs",3,1,105.8801323333333,317.640397,0.8234679448457706,0.8234679448457706,3,1,17.37449166666667,52.123475,0.1547081917640233,0.1547081917640233,1839.613476891064,16556.52129201957,0.1273972367227257,0.1273972367227257),
    @("M2","Tnc","Egfr","ECs",2,0.6666666666666666,0.1148696666666667,0.344609,0.0008933827928862465,0.0008933827928862465,3,1,1.854571666666667,5.563715,0.01651371646154392,0.01651371646154392,0.2130340291594444,1.917306262435,0.00001475307013334569,0.00001475307013334569),
    @("M2","Tnc","Egfr","FAPs",2,0.6666666666666666,0.1148696666666667,0.344609,0.0008933827928862465,0.0008933827928862465,3,1,92.91372433333333,278.741173,0.8273343794712995,0.8273343794712996,10.67296854292855,96.056716886357,0.0007391262985828792,0.0007391262985828793),
    @("M2","Tnc","Egfr","M2",2,0.6666666666666666,0.1148696666666667,0.344609,0.0008933827928862465,0.0008933827928862465,3,1,0.162136,0.4864080000000001,0.001443712303133186,0.001443712303133187,0.01862450827466667,0.167620574472,0.000001289787729497361,0.000001289787729497362),
    @("M2","Tnc","Egfr","sCs",2,0.6666666666666666,0.1148696666666667,0.344609,0.0008933827928862465,0.0008933827928862465,3,1,17.37449166666667,52.123475,0.1547081917640233,0.1547081917640233,1.995802066252778,17.962218596275,0.0001382136364405242,0.0001382136364405242),
    @("sCs","Tnc","Egfr","ECs",3,1,13.635153,40.905459,0.1060454985380935,0.1060454985380935,3,1,1.854571666666667,5.563715,0.01651371646154392,0.01651371646154392,25.287368424465,227.586315820185,0.001751205294881146,0.001751205294881146),
    @("sCs","Tnc","Egfr","FAPs",3,1,13.635153,40.905459,0.1060454985380935,0.1060454985380935,3,1,92.91372433333333,278.741173,0.8273343794712995,0.8273343794712996,1266.892847084823,11402.03562376341,0.08773508672873816,0.08773508672873817),
    @("sCs","Tnc","Egfr","M2",3,1,13.635153,40.905459,0.1060454985380935,0.1060454985380935,3,1,0.162136,0.4864080000000001,0.001443712303133186,0.001443712303133187,2.210749166808001,19.896742501272,0.0001530991909313379,0.0001530991909313379),
    @("sCs","Tnc","Egfr","sCs",3,1,13.635153,40.905459,0.1060454985380935,0.1060454985380935,3,1,17.37449166666667,52.123475,0.1547081917640233,0.1547081917640233,236.903852172225,2132.134669550025,0.01640610732354282,0.01640610732354282)
)

# Clear any previously existing data rows below the header before writing
# the refreshed table (the new table is longer: 16 rows vs. 12 before).
$ws.Range("A2:T13").Clear() | Out-Null

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowValues[$c]
    }
}
